$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("3:3").Delete() | Out-Null
